# [Refactor] format tables for submission
#
# Adds the table title to row 1 and the data-source footnote to row 29
# (both rows were previously blank / outside the sheet's used range),
# then leaves the selection where Excel would land after typing the
# footnote and pressing Enter (one row below it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Table 3: Age- and sex-adjusted prevalence of risk factors for survey respondents aged 12 years and older, from 2000 to 2018."
$ws.Range("A29").Value = "Data source: CCHS"

[void]$ws.Range("A30").Select()
